$wb = $excel.ActiveWorkbook
$wsInv = $wb.Worksheets.Item("Investment_Cost")
$wsSrc = $wb.Worksheets.Item("Sources")

# --- Correct the investment-cost header years (the actual purpose of the commit) ---
# Column B keeps "...Value 2020"; C/D/E/F were wrongly duplicating 2020 and are
# corrected to 2025 / 2030 / 2040 / 2050 respectively.
$wsInv.Range("C1").Value = "Investment_Cost [Euro/MW or MWh] Value 2025"
$wsInv.Range("D1").Value = "Investment_Cost [Euro/MW or MWh] Value 2030"
$wsInv.Range("E1").Value = "Investment_Cost [Euro/MW or MWh] Value 2040"
$wsInv.Range("F1").Value = "Investment_Cost [Euro/MW or MWh] Value 2050"

# --- Update the selected / active cells to match the state the file was saved in ---
$wsInv.Activate() | Out-Null
$wsInv.Range("G7").Select() | Out-Null

# --- S&P Capital IQ Excel add-in workbook-level defined names (injected by the
#     add-in when the file was opened/saved with it active) ---
$n = $wb.Names.Add("CIQWBGuid", '="2f1f32c2-b0bb-4653-88c7-8313657ec4b2"')
$n.Visible = $false
$n = $wb.Names.Add("CIQWBInfo", '="{ ""CIQVersion"":""9.51.3510.3078"" }"')
$n.Visible = $false

$wb.Names.Add("IQ_CH", "=110000")
$wb.Names.Add("IQ_CQ", "=5000")
$wb.Names.Add("IQ_CY", "=10000")
$wb.Names.Add("IQ_DAILY", "=500000")
$n = $wb.Names.Add("IQ_DNTM", "=700000")
$n.Visible = $false
$wb.Names.Add("IQ_FH", "=100000")
$wb.Names.Add("IQ_FQ", "=500")
$n = $wb.Names.Add("IQ_FWD_CY", "=10001")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_CY1", "=10002")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_CY2", "=10003")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_FY", "=1001")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_FY1", "=1002")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_FY2", "=1003")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_Q", "=501")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_Q1", "=502")
$n.Visible = $false
$n = $wb.Names.Add("IQ_FWD_Q2", "=503")
$n.Visible = $false
$wb.Names.Add("IQ_FY", "=1000")
$n = $wb.Names.Add("IQ_LATESTK", "=1000")
$n.Visible = $false
$n = $wb.Names.Add("IQ_LATESTQ", "=500")
$n.Visible = $false
$wb.Names.Add("IQ_LTM", "=2000")
$n = $wb.Names.Add("IQ_LTMMONTH", "=120000")
$n.Visible = $false
$wb.Names.Add("IQ_MONTH", "=15000")
$n = $wb.Names.Add("IQ_MTD", "=800000")
$n.Visible = $false
$n = $wb.Names.Add("IQ_NAMES_REVISION_DATE_", "=45429.7589699074")
$n.Visible = $false
$wb.Names.Add("IQ_NTM", "=6000")
$n = $wb.Names.Add("IQ_QTD", "=750000")
$n.Visible = $false
$n = $wb.Names.Add("IQ_TODAY", "=0")
$n.Visible = $false
$wb.Names.Add("IQ_WEEK", "=50000")
$wb.Names.Add("IQ_YTD", "=3000")
$n = $wb.Names.Add("IQ_YTDMONTH", "=130000")
$n.Visible = $false

# --- Sources sheet no longer keeps an explicit A2 selection ---
$wsSrc.Activate() | Out-Null
$wsSrc.Range("A1").Select() | Out-Null
$wsInv.Activate() | Out-Null
